# Reading Bordalo et al and Hebden et al: add two new rows (13 and 14)
# of notes to the "Record of empirical features of expectations" sheet,
# and move the current selection to C15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: diagnostic expectations (Bordalo et al.) ---
# Write B13/C13 before A13 so that new shared-string entries are created
# in the same order as in the target workbook (41, 42, 43, ...).
$ws.Range("B13").Value = "explain that individuals overreact to their news using representativeness heuristic. When everyone overreacts, the average underreacts to news."
$ws.Range("C13").Value = "cannot account for time-varying over and underreaction, silent on LR-E, both of which is to say that it cannot account for anchoring"
$ws.Range("A13").Value = "diagnostic expectations (Bordalo et al, 2018 (unpublished) and 2018)"

# --- Row 14: Hebden et al 2020 ---
$ws.Range("A14").Value = "Hebden et al 2020"
$ws.Range("B14").Value = "Makeup strategies work well when expectations understand, believe and act on the CB's policy commitment. Otherwise, costly action is necessary (as in Goodfriend and King)"
$ws.Range("C14").Value = "Long-run expectations seem anchored generally, but in surveys, a majority of respondents does revise its expectations for persistent deviations . In survey data, inflation expectations, even individual ones, generally underreact to news."

# Row heights grew to fit the new wrapped text.
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 46

# Selection moved to C15 (just below the newly added notes).
$ws.Range("C15").Select() | Out-Null
